# tbidbaxlipo/data - "Updated dataset from Justin" correction
#
# Justin's 2/13/2015 email: the 126C FRET/NBD values originally sent on
# 12/3 were accidentally a moving average of 5 points rather than the
# absolute per-timepoint calculations. This replaces the FRET (col C) and
# NBD (col D) values for the Bax 126C series with the corrected data, fixes
# up the header labels, and trims the now-unsupported trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 1 merged group headers: "Bax 126C"/"Bax 54C" -> "126C"/"54C"
$ws.Cells.Item(1, 2).Value = "126C"
$ws.Cells.Item(1, 6).Value = "54C"

# Row 2 column headers - relabel the X/TB/duplicate-NBD columns and add the
# new "Time" label over column A
$ws.Cells.Item(2, 1).Value = "Time"
$ws.Cells.Item(2, 2).Value = "RLS"
$ws.Cells.Item(2, 3).Value = "FRET"
$ws.Cells.Item(2, 4).Value = "NBD"
$ws.Cells.Item(2, 5).Value = "Time"
$ws.Cells.Item(2, 6).Value = "RLS"
$ws.Cells.Item(2, 7).Value = "DAC"
$ws.Cells.Item(2, 8).Value = "NBD"

# Corrected FRET (col C) and NBD (col D) values for the Bax 126C series,
# row -> (FRET, NBD)
$newData = @(
    @(3,0,1),
    @(4,6.5132786519519259,1.6024160797200757),
    @(5,6.3177301462430835,1.7715711395971492),
    @(6,10.037134755078526,1.9107581630316364),
    @(7,13.957034394680679,2.0671207980043342),
    @(8,11.321493945916783,2.2103309883308917),
    @(9,14.164013744442649,2.3338674459841156),
    @(10,18.188425085002812,2.4743174638896219),
    @(11,18.716556932224691,2.5523828024526649),
    @(12,19.05159528954341,2.6690050373555243),
    @(13,20.83992117224096,2.7159270565909925),
    @(14,22.672274725265428,2.8116807147485279),
    @(15,24.967186484181269,2.8970634721693806),
    @(16,27.166672281750216,2.9957847560658362),
    @(17,26.590211568086154,3.0932638436423265),
    @(18,25.672712316695943,3.1344895079828765),
    @(19,24.448324391310383,3.1785236224582563),
    @(20,26.226156084266261,3.2568329273510472),
    @(21,24.404985756844333,3.348221095586736),
    @(22,31.750133657988577,3.3956523467386308),
    @(23,25.352282320910479,3.437990375832793),
    @(24,25.526663556548112,3.5466996523747283),
    @(25,39.46487113384277,3.5468373436287268),
    @(26,28.69992477360266,3.652009733643041),
    @(27,29.37910297209071,3.6550929162025705),
    @(28,26.427608827374726,3.6976638918116933),
    @(29,24.432049097160057,3.7237168828551028),
    @(30,36.541679742388332,3.793144424694904),
    @(31,31.148375751960511,3.8394718885811407),
    @(32,34.912258439868694,3.8860393866960683),
    @(33,34.821350622954427,3.9004123221031821),
    @(34,34.240654617605728,3.9321095259436847),
    @(35,34.451911249113657,4.0010941082476723),
    @(36,34.932981364404228,3.9976835636031458),
    @(37,33.645259759078719,4.0393972418636039),
    @(38,36.556525170952192,3.9926601352105542),
    @(39,36.725965876933131,4.0473261362942523),
    @(40,35.050113679316667,4.1395418966834807),
    @(41,32.194971757183012,4.1849718844616888),
    @(42,35.643540671252651,4.1846029621903202),
    @(43,33.046671697365845,4.2243676546089386),
    @(44,39.821756879590588,4.1833576012287477),
    @(45,39.184222743281282,4.1917527498467058),
    @(46,42.71772129558736,4.2569409359822581),
    @(47,39.056252393387922,4.2213686939521828),
    @(48,39.628544021353044,4.2580942469546814),
    @(49,37.631836763511338,4.2970410533245493),
    @(50,45.081038420655204,4.3638756205558602),
    @(51,41.21841569474973,4.3485038592488401),
    @(52,38.922534987563317,4.3845113410755703),
    @(53,37.373455988362615,4.4785276064756872),
    @(54,37.162248339104984,4.3939296485724189),
    @(55,40.058986504833861,4.4415812508753838),
    @(56,41.493807181717699,4.4356704587816154),
    @(57,41.305944438803706,4.4586681486158142),
    @(58,39.790993205938854,4.4394573294239557),
    @(59,43.32449493478515,4.4725695479088259),
    @(60,42.612756363383056,4.4702604881516113),
    @(61,41.907053190128806,4.488735358877018),
    @(62,39.918786434609707,4.4429204224683687),
    @(63,42.476523344980322,4.5271212426793683),
    @(64,45.0481950677339,4.5599691391012351),
    @(65,51.473598900676208,4.5064721593363188),
    @(66,40.767078412241197,4.6107986555464802),
    @(67,39.779527590762562,4.698991619229961),
    @(68,42.511778610625107,4.5961029357242564),
    @(69,39.255308914354494,4.6450454951089704),
    @(70,41.46443700447464,4.648974120025918),
    @(71,48.458889252343226,4.63127212610984),
    @(72,42.023271514480989,4.6731494538115257),
    @(73,42.040602375979752,4.6574331032121226),
    @(74,55.772237613025212,4.6773098063372514),
    @(75,47.111709036578873,4.6253101851010623),
    @(76,44.561241292552403,4.6615930945804616),
    @(77,46.825133820873766,4.6755992295158588),
    @(78,40.896634874533746,4.6464784773333676),
    @(79,46.35208863889283,4.6837817015772307),
    @(80,48.522397146253091,4.6654297146529347),
    @(81,46.239088982486265,4.6836893355753686),
    @(82,47.337655711453444,4.7123526367818718),
    @(83,44.632854497175799,4.7260847435494755),
    @(84,45.660763226911605,4.7160978840352978),
    @(85,45.09678111543758,4.7154041458385132),
    @(86,43.640632609917596,4.7102721445062938),
    @(87,50.865930289252404,4.7170688556716893),
    @(88,45.790721083761184,4.7854150922932064),
    @(89,48.212810677062102,4.8050720762926922),
    @(90,43.453612295040621,4.7587768908479671),
    @(91,45.652849355442939,4.760487828826748),
    @(92,46.25912634889113,4.7560022540735458),
    @(93,47.041854246616452,4.7708470867521493),
    @(94,44.441810238015954,4.8072918850305122),
    @(95,47.535103974013182,4.7846289429498867),
    @(96,52.97258414303991,4.7600254119366809),
    @(97,46.131205125527153,4.755771248779542),
    @(98,45.364073370265302,4.7007945620532166),
    @(99,46.946536994838354,4.728905247312114),
    @(100,44.599933922806045,4.7000545957108271),
    @(101,46.365043756996485,4.749528237033914),
    @(102,46.372856446776382,4.7595168571903566),
    @(103,48.435025325258117,4.8243604540335241),
    @(104,46.39049528096465,4.7682110892698679),
    @(105,46.193782850253484,4.792074969953398),
    @(106,50.611408433077322,4.7810677956747574),
    @(107,47.072539966918171,4.7627076827092401),
    @(108,47.968413143077981,4.7696444326516527),
    @(109,42.051376550564136,4.7884676396773385),
    @(110,54.567401171499696,4.9246687582946285),
    @(111,43.405834802188039,4.8455474361606354),
    @(112,46.860000816760405,4.829032611722063),
    @(113,57.217975331711123,4.8003078686150058),
    @(114,47.867218253568275,4.8452232974053215),
    @(115,48.522628152448036,4.7899475272174419),
    @(116,48.146090146048806,4.7954975230810613),
    @(117,48.548285907761382,4.8285698788192821),
    @(118,48.950428545518818,4.8026666778023532),
    @(119,50.2387618645215,4.7649737195926649),
    @(120,48.819069149223083,4.7800964628809774),
    @(121,48.335079115445524,4.763863792651418),
    @(122,51.733297436284118,4.806320597381406),
    @(123,49.13626247465708,4.7966539941806259),
    @(124,48.292627062363913,4.7704305819949733),
    @(125,53.899485268691571,4.8186246426940178),
    @(126,52.092577254692209,4.8037769659009868),
    @(127,48.039183154789157,4.7978563322684087),
    @(128,49.353707327521221,4.8280147347699645)
)

foreach ($row in $newData) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
}

# Remove the now-obsolete trailing rows (129-174); the corrected dataset
# only covers time points through row 128
$ws.Rows("129:174").Delete()

# Restore the cursor/selection position recorded in the corrected workbook
[void]$ws.Range("G23").Select()
